# Insert a new "chemical_recycling_pyrolysis" parameter row right after the
# existing "chemical_recycling_gasification" row (currently row 9), pushing
# every row from the old row 10 onward down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a full new row at row 10 — shifts rows 10:24 down to 11:25.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the pyrolysis parameter + its value.
$ws.Cells.Item(10, 1).Value = "chemical_recycling_pyrolysis"
$ws.Cells.Item(10, 2).Value = $true
